# Updated cryptos list on Tue Feb 20 00:48:23 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.583.66"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "2.935.14"
$ws.Range("E3").Value = "  +2.02%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'351.11"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").Value = "'111.11"
$ws.Range("E6").Value = "  -1.34%  "

$ws.Range("D7").Value = "'0.561"
$ws.Range("E7").Value = "  +0.56%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.626"
$ws.Range("E9").Value = "  +0.27%  "

$ws.Range("D10").Value = "'39.09"
$ws.Range("E10").Value = "  -3.30%  "

$ws.Range("D11").Value = "'0.0903"
$ws.Range("E11").Value = "  +5.31%  "

$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("D13").Value = "'19.71"
$ws.Range("E13").Value = "  -2.20%  "

$ws.Range("D14").Value = "'8.02"
$ws.Range("E14").Value = "  +1.92%  "

$ws.Range("D15").Value = "3.416.66"
$ws.Range("E15").Value = "  +2.27%  "

$ws.Range("D16").Value = "2.958.72"
$ws.Range("E16").Value = "  +2.66%  "

$ws.Range("D17").Value = "'0.989"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "51.812.28"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").Value = "'7.61"
$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").Value = "'14.36"
$ws.Range("E20").Value = "  +4.90%  "

$ws.Range("D21").Value = "'3.23"
$ws.Range("E21").Value = "  -3.53%  "

$ws.Range("D22").Value = "0.0₃0983"
$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("D23").Value = "'70.93"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "'270.80"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "'2.77"
$ws.Range("E25").Value = "  -1.23%  "

$ws.Range("D26").Value = "'0.180"
$ws.Range("E26").Value = "  +9.32%  "

$ws.Range("D27").Value = "'27.09"
$ws.Range("E27").Value = "  +2.38%  "

$ws.Range("E28").Value = "  +0.19%  "

$ws.Range("D29").Value = "'7.51"
$ws.Range("E29").Value = "  +19.06%  "

$ws.Range("D30").Value = "'0.109"
$ws.Range("E30").Value = "  +18.87%  "

$ws.Range("D31").Value = "'10.63"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'2.26"
$ws.Range("E32").Value = "  +9.05%  "

$ws.Range("D33").Value = "'37.08"
$ws.Range("E33").Value = "  -6.80%  "

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'6.21"
$ws.Range("E34").Value = "  +4.69%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'52.57"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.0441"
$ws.Range("E36").Value = "  -5.50%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").Value = "'3.36"
$ws.Range("E38").Value = "  +1.50%  "

$ws.Range("D39").Value = "'18.61"
$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").Value = "'2.03"
$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("E42").Value = "  +1.10%  "

$ws.Range("D43").Value = "'23.51"
$ws.Range("E43").Value = "  +4.29%  "

$ws.Range("E44").Value = "  -1.71%  "

$ws.Range("E45").Value = "  +2.19%  "

$ws.Range("D46").Value = "'3.48"
$ws.Range("E46").Value = "  -1.29%  "

$ws.Range("D47").Value = "2.150.93"
$ws.Range("E47").Value = "  -1.85%  "

$ws.Range("D48").Value = "'113.48"
$ws.Range("E48").Value = "  -7.23%  "

$ws.Range("D49").Value = "'0.243"
$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").Value = "'0.0330"
$ws.Range("E50").Value = "  +2.80%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'9.15"
$ws.Range("E51").Value = "  +3.45%  "
